# SDM_Avaliacao.xlsx — weekly report history update
#
# The evaluator "Davan Martinho do Nascimento" (row 20) is removed from the
# roster entirely (he no longer has an evaluation row). Deleting the whole
# worksheet row shifts every subsequent evaluator up by one row, drops the
# now-unused shared string, and re-derives the summary row's SUM ranges
# automatically. The only thing Excel's row-shift can't infer on its own is
# that the hard-coded "/28" divisor (count of evaluators) in the summary
# formulas must become "/27" now that there are 27 evaluators instead of 28
# — except for the handful of columns whose column total is zero, where the
# author's original formula kept the stale "/28" divisor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove Davan Martinho do Nascimento's row; everything below shifts up.
$ws.Rows("20:20").Delete()

# Summary row is now row 29 (was row 30). Re-point each average formula at
# the new 27-row data range (D2:D28 etc.) and fix the literal evaluator
# count used as the divisor.
$ws.Range("D29").Formula = "=SUM(D2:D28)*100/27"
$ws.Range("E29").Formula = "=SUM(E2:E28)*100/27"
$ws.Range("F29").Formula = "=SUM(F2:F28)*100/27"
$ws.Range("G29").Formula = "=SUM(G2:G28)*100/27"
$ws.Range("H29").Formula = "=SUM(H2:H28)*100/28"
$ws.Range("I29").Formula = "=SUM(I2:I28)*100/27"
$ws.Range("K29").Formula = "=SUM(K2:K28)*100/27"
$ws.Range("L29").Formula = "=SUM(L2:L28)*100/27"
$ws.Range("M29").Formula = "=SUM(M2:M28)*100/28"
$ws.Range("N29").Formula = "=SUM(N2:N28)*100/27"
$ws.Range("O29").Formula = "=SUM(O2:O28)*100/27"
$ws.Range("Q29").Formula = "=SUM(Q2:Q28)*100/28"
$ws.Range("R29").Formula = "=SUM(R2:R28)*100/28"
$ws.Range("S29").Formula = "=SUM(S2:S28)*100/27"
$ws.Range("T29").Formula = "=SUM(T2:T28)*100/27"
$ws.Range("U29").Formula = "=SUM(U2:U28)*100/27"
$ws.Range("V29").Formula = "=SUM(V2:V28)*100/28"
$ws.Range("W29").Formula = "=SUM(W2:W28)*100/27"
$ws.Range("X29").Formula = "=SUM(X2:X28)*100/27"
$ws.Range("Y29").Formula = "=SUM(Y2:Y28)*100/27"
$ws.Range("Z29").Formula = "=SUM(Z2:Z28)*100/28"
$ws.Range("AA29").Formula = "=SUM(AA2:AA28)*100/27"
$ws.Range("AB29").Formula = "=SUM(AB2:AB28)*100/27"
$ws.Range("AC29").Formula = "=SUM(AC2:AC28)*100/27"
$ws.Range("AD29").Formula = "=SUM(AD2:AD28)*100/27"
$ws.Range("AE29").Formula = "=SUM(AE2:AE28)*100/27"
$ws.Range("AF29").Formula = "=SUM(AF2:AF28)*100/27"
$ws.Range("AG29").Formula = "=SUM(AG2:AG28)*100/27"
$ws.Range("AH29").Formula = "=SUM(AH2:AH28)*100/27"
$ws.Range("AI29").Formula = "=SUM(AI2:AI28)*100/27"

# Restore the view: scrolled slightly right and selection left on the
# summary row, matching where the author's edit left the cursor.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("F30").Select()
